$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append newly logged/simulated per-play yardage values to the
# running lists stored in B2 (R/OFF), C2 (R/DEF), B3 (P/OFF), C3 (P/DEF).
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 6 3 15 15 4 1 0 8 3 10 3 0 4 1 0 4 4 2 3"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 2 4 9 8 4 0 8 -2 0 13 0 11 4 9 3 -2 1 6 3 4 0 5 2 7 1 8 2"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 14 48 8 2 21 -4 32 15 -2 19 5 5 6 12"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 7 7 5 10 24 6 5 7 14 22 4 10 9 5 40 8 9 25 16 13 1 18 11 12 8 13 8 61 9"

# ---------------------------------------------------------------------------
# OFF sheet: updated rush/pass attempt-by-situation tallies and totals.
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 181
$wsOFF.Range("E2").Value = 6
$wsOFF.Range("F2").Value = 72
$wsOFF.Range("J2").Value = 40
$wsOFF.Range("N2").Value = 22
$wsOFF.Range("O2").Value = 30
$wsOFF.Range("P2").Value = 15

$wsOFF.Range("C3").Value = 178
$wsOFF.Range("F3").Value = 87
$wsOFF.Range("G3").Value = 31
$wsOFF.Range("H3").Value = 26
$wsOFF.Range("I3").Value = 54
$wsOFF.Range("J3").Value = 49
$wsOFF.Range("L3").Value = 247
$wsOFF.Range("M3").Value = 160
$wsOFF.Range("Q3").Value = 490

# ---------------------------------------------------------------------------
# DEF sheet: updated rush/pass attempt-by-situation tallies and totals.
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 169
$wsDEF.Range("F2").Value = 62
$wsDEF.Range("G2").Value = 37
$wsDEF.Range("H2").Value = 5
$wsDEF.Range("I2").Value = 5
$wsDEF.Range("J2").Value = 30
$wsDEF.Range("N2").Value = 18
$wsDEF.Range("O2").Value = 14
$wsDEF.Range("P2").Value = 8

$wsDEF.Range("C3").Value = 185
$wsDEF.Range("D3").Value = 2
$wsDEF.Range("E3").Value = 27
$wsDEF.Range("F3").Value = 112
$wsDEF.Range("G3").Value = 37
$wsDEF.Range("H3").Value = 20
$wsDEF.Range("I3").Value = 50
$wsDEF.Range("J3").Value = 68
$wsDEF.Range("L3").Value = 280
$wsDEF.Range("M3").Value = 202
$wsDEF.Range("Q3").Value = 498

# ---------------------------------------------------------------------------
# ST sheet: updated kickoff/punt counters, and appended per-kick/punt
# distance and return logs (D / RA / RM rows for KO and PT).
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 67
$wsST.Range("D2").Value = 54
$wsST.Range("F2").Value = 98
$wsST.Range("G2").Value = 86
$wsST.Range("B3").Value = 43

$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 63"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 20"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 0"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 50 56 42 53 53 54"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 21 0 0 0 0 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 7 7 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: updated road interception/fumble counts.
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value = 8
$wsTURNS.Range("D3").Value = 10
